# Regenerate the task-order sheets (new experiment order generation run).
#
# Original tabs (left to right): GNG_TO-..., NB_TO-..., RS_TO-..., TOL_TO-..., vSAT_TO-...
# New tabs      (left to right): NB_TO-...,  TOL_TO-..., GNG_TO-..., RS_TO-..., vSAT_TO-...
#
# Each tab keeps its task-family prefix but gets a brand new numeric suffix
# (new generation timestamp) and brand new generated stim-file rows.

$wb = $excel.ActiveWorkbook

# ---- capture the five original sheets by their original (pre-edit) names ----
$sheetGNG  = $wb.Worksheets.Item("GNG_TO-16512556056162512")
$sheetNB   = $wb.Worksheets.Item("NB_TO-16512556077562475")
$sheetRS   = $wb.Worksheets.Item("RS_TO-1651255607758249")
$sheetTOL  = $wb.Worksheets.Item("TOL_TO-16512556078212523")
$sheetVSAT = $wb.Worksheets.Item("vSAT_TO-16512556078982496")

# ---- rename every tab to its freshly generated name ----
$sheetNB.Name   = "NB_TO-16515890331842973"
$sheetTOL.Name  = "TOL_TO-1651589033231173"
$sheetGNG.Name  = "GNG_TO-16515890332624226"
$sheetRS.Name   = "RS_TO-16515890332624226"
$sheetVSAT.Name = "vSAT_TO-16515890333249228"

# =========================================================================
# NB_TO-16515890331842973  -> new-balance task-order rows (9 rows, A1:B10)
# =========================================================================
$ws = $wb.Worksheets.Item("NB_TO-16515890331842973")
$nbFiles = @(
    "ZB-match_4-1651589031995488.csv",
    "TB-1651589033054729.csv",
    "TB-16515890331686702.csv",
    "ZB-match_2-165158903177121.csv",
    "OB-16515890327214031.csv",
    "TB-16515890329578586.csv",
    "ZB-match_4-16515890319798965.csv",
    "OB-16515890322477357.csv",
    "OB-16515890322816722.csv"
)
# extend the styled index column (A) down to row 10 (it already covers rows 2-5)
$ws.Range("A2").Copy()
$ws.Range("A6:A10").PasteSpecial(-4122)
for ($i = 0; $i -lt $nbFiles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $nbFiles[$i]
}

# =========================================================================
# TOL_TO-1651589033231173 -> tolerance task-order rows (6 rows, A1:B7)
# =========================================================================
$ws = $wb.Worksheets.Item("TOL_TO-1651589033231173")
$tolFiles = @(
    "MM_stims-1651589033199923.csv",
    "ZM_stims-16515890331842973.csv",
    "MM_stims-16515890332155497.csv",
    "ZM_stims-1651589033199923.csv",
    "MM_stims-1651589033231173.csv",
    "ZM_stims-16515890332155497.csv"
)
for ($i = 0; $i -lt $tolFiles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $tolFiles[$i]
}
# drop the trailing three rows (sheet used to have 8 data rows, now only 6)
$ws.Range("A8:B10").EntireRow.Delete()

# =========================================================================
# GNG_TO-16515890332624226 -> go/no-go task-order rows (4 rows, A1:B5)
# =========================================================================
$ws = $wb.Worksheets.Item("GNG_TO-16515890332624226")
$gngFiles = @(
    "go_stims-1651589033231173.csv",
    "GNG_stims-16515890332467983.csv",
    "go_stims-16515890332467983.csv",
    "GNG_stims-16515890332624226.csv"
)
for ($i = 0; $i -lt $gngFiles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $gngFiles[$i]
}

# =========================================================================
# RS_TO-16515890332624226 -> resting-state task-order rows (swap the two)
# =========================================================================
$ws = $wb.Worksheets.Item("RS_TO-16515890332624226")
$ws.Cells.Item(2, 2).Value = "eyes closed"
$ws.Cells.Item(3, 2).Value = "eyes open"

# =========================================================================
# vSAT_TO-16515890333249228 -> vSAT task-order rows (4 rows, A1:B5)
# =========================================================================
$ws = $wb.Worksheets.Item("vSAT_TO-16515890333249228")
$vsatFiles = @(
    "vSAT_stims-16515890333092985.csv",
    "vSAT_stims-1651589033293672.csv",
    "SAT_stims-16515890332780483.csv",
    "SAT_stims-16515890332624226.csv"
)
for ($i = 0; $i -lt $vsatFiles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $vsatFiles[$i]
}

# =========================================================================
# Reorder tabs: NB, TOL, GNG, RS, vSAT (left to right)
# =========================================================================
$target = $wb.Worksheets.Item("GNG_TO-16515890332624226")
$wb.Worksheets.Item("NB_TO-16515890331842973").Move($target)

$target = $wb.Worksheets.Item("GNG_TO-16515890332624226")
$wb.Worksheets.Item("TOL_TO-1651589033231173").Move($target)
